$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "64.030.55"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "3.170.26"
$ws.Range("E3").Value = "  -8.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'563.52"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.609"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "3.162.80"
$ws.Range("E9").Value = "  -8.50%  "
$ws.Range("E10").Value = "  -6.73%  "
$ws.Range("D11").Value = "'6.62"
$ws.Range("E11").Value = "  -4.83%  "
$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  -5.63%  "
$ws.Range("D13").Value = "3.714.07"
$ws.Range("E13").Value = "  -8.46%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "'27.30"
$ws.Range("E15").Value = "  -7.23%  "
$ws.Range("D16").Value = "64.021.38"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("E17").Value = "  -5.78%  "
$ws.Range("D18").Value = "3.168.96"
$ws.Range("E18").Value = "  -8.27%  "
$ws.Range("D19").Value = "'5.71"
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("D20").Value = "'13.05"
$ws.Range("E20").Value = "  -5.34%  "
$ws.Range("D21").Value = "'353.73"
$ws.Range("E21").Value = "  -5.05%  "
$ws.Range("D22").Value = "'7.21"
$ws.Range("E22").Value = "  -5.61%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'69.03"
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").Value = "'0.502"
$ws.Range("E25").Value = "  -6.76%  "
$ws.Range("D26").Value = "'0.0000117"
$ws.Range("E26").Value = "  -4.73%  "
$ws.Range("D27").Value = "'9.59"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "'5.63"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("D33").Value = "'22.10"
$ws.Range("E33").Value = "  -6.81%  "
$ws.Range("D34").Value = "'6.62"
$ws.Range("E34").Value = "  -6.18%  "
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.44"
$ws.Range("E36").Value = "  -7.83%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'155.14"
$ws.Range("E37").Value = "  -4.14%  "
$ws.Range("D38").Value = "'0.808"
$ws.Range("E38").Value = "  -8.43%  "
$ws.Range("D39").Value = "'25.89"
$ws.Range("E39").Value = "  -9.05%  "
$ws.Range("D40").Value = "'2.54"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").Value = "2.603.18"
$ws.Range("E42").Value = "  -6.63%  "
$ws.Range("D43").Value = "'4.18"
$ws.Range("E43").Value = "  -7.08%  "
$ws.Range("E44").Value = "  -5.83%  "
$ws.Range("D45").Value = "'0.0658"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'39.12"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'329.31"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "'23.91"
$ws.Range("E48").Value = "  -5.43%  "
$ws.Range("E49").Value = "  -7.56%  "
$ws.Range("D51").Value = "'0.998"
